# Apply the scheduled cryptos.xlsx data refresh (GitHub Actions bot commit:
# 'Updated cryptos list ... with GitHub Actions'): refreshed Price (column D)
# and Volume(1h) (column E) figures for every listed coin, plus two adjacent-
# rank swaps where the refresh reordered rows:
#   rows 10-11: Cardano <-> Dogecoin
#   rows 50-51: Monero <-> THORChain
#
# Every Price/Volume cell in this sheet is stored as literal text (e.g.
# '0.997', '0.0000280') even though many look numeric, so any replacement
# that parses as a plain number is written with a leading apostrophe to stop
# Excel from re-parsing it into a Number cell (which would also silently drop
# meaningful trailing zeros, e.g. '137.10' -> 137.1). The style is restored to
# Normal right after so the apostrophe's quote-prefix formatting doesn't stick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.250.37'
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('D3').Value = '3.404.14'
$ws.Range('E3').Value = '  +2.82%  '
$ws.Range('D4').Value = '''0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '''566.57'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.45%  '
$ws.Range('D6').Value = '''178.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.64%  '
$ws.Range('E7').Value = '  +4.50%  '
$ws.Range('D8').Value = '3.389.78'
$ws.Range('E8').Value = '  +2.61%  '
$ws.Range('D9').Value = '''0.997'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '''0.167'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.33%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').Value = '''0.638'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.85%  '
$ws.Range('D12').Value = '''55.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.50%  '
$ws.Range('D13').Value = '''0.0000280'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.97%  '
$ws.Range('D14').Value = '''9.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.66%  '
$ws.Range('D15').Value = '3.903.70'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '''18.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.54%  '
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '3.357.41'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '''11.93'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.40%  '
$ws.Range('D20').Value = '64.846.79'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').Value = '''0.998'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.06%  '
$ws.Range('D22').Value = '''469.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +18.36%  '
$ws.Range('D23').Value = '''4.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +15.26%  '
$ws.Range('D24').Value = '''4.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.08%  '
$ws.Range('D25').Value = '''86.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.76%  '
$ws.Range('D26').Value = '''13.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('D27').Value = '''10.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.87%  '
$ws.Range('D28').Value = '''2.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.24%  '
$ws.Range('D29').Value = '''8.91'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.76%  '
$ws.Range('D30').Value = '''30.59'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.31%  '
$ws.Range('D31').Value = '''6.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.18%  '
$ws.Range('D32').Value = '''11.58'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.71%  '
$ws.Range('D33').Value = '''580.94'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').Value = '''0.110'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.09%  '
$ws.Range('D35').Value = '''60.51'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.69%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '''0.142'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').Value = '''36.25'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('D39').Value = '0.0₃0769'
$ws.Range('E39').Value = '  +6.06%  '
$ws.Range('D40').Value = '''3.49'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('D41').Value = '''0.375'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.80%  '
$ws.Range('D42').Value = '3.107.84'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = '''0.997'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = '''2.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.37%  '
$ws.Range('D45').Value = '''2.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.18%  '
$ws.Range('D46').Value = '''0.0414'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.82%  '
$ws.Range('D47').Value = '''3.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('E48').Value = '  +5.17%  '
$ws.Range('D49').Value = '''2.57'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''8.46'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.16%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''137.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.29%  '
